$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume/Number and date range) ---
$ws.Range("A8").Value = "Volume 31   Number  52"
$ws.Range("C9").Value = "Report Covering the Week  12/23/2024  Through  12/29/2024"

# --- Column width updates (I, J now match H) ---
$ws.Columns("I").ColumnWidth = $ws.Columns("H").ColumnWidth
$ws.Columns("J").ColumnWidth = $ws.Columns("H").ColumnWidth

# --- Reference cells used to clone the "0" / "***.*" placeholder text+style ---
$srcZero = $ws.Range("C14")
$srcStar = $ws.Range("E14")

# --- Cells switching from a numeric value to the "0" placeholder text ---
$srcZero.Copy($ws.Range("G15"))
$srcZero.Copy($ws.Range("D22"))
$srcZero.Copy($ws.Range("C25"))
$srcZero.Copy($ws.Range("G27"))
$srcZero.Copy($ws.Range("C28"))

# --- Cells switching from a numeric value to the "***.*" placeholder text ---
$srcStar.Copy($ws.Range("H15"))
$srcStar.Copy($ws.Range("E22"))
$srcStar.Copy($ws.Range("H27"))

# --- Updated crime-statistics figures ---
$ws.Range("M14").Value = 40
$ws.Range("F15").Value = 3
$ws.Range("I15").Value = 21
$ws.Range("K15").Value = 50
$ws.Range("L15").Value = 16.666666666666
$ws.Range("M15").Value = -22.222222222222
$ws.Range("N15").Value = -66.666666666666
$ws.Range("C16").Value = 3
$ws.Range("D16").Value = 5
$ws.Range("E16").Value = -40
$ws.Range("G16").Value = 21
$ws.Range("H16").Value = 19.047619047619
$ws.Range("I16").Value = 248
$ws.Range("J16").Value = 205
$ws.Range("K16").Value = 20.975609756097
$ws.Range("L16").Value = 31.914893617021
$ws.Range("M16").Value = -1.587301587301
$ws.Range("N16").Value = -73.131094257854
$ws.Range("C17").Value = 11
$ws.Range("D17").Value = 5
$ws.Range("E17").Value = 120
$ws.Range("F17").Value = 29
$ws.Range("G17").Value = 23
$ws.Range("H17").Value = 26.086956521739
$ws.Range("I17").Value = 482
$ws.Range("J17").Value = 428
$ws.Range("K17").Value = 12.616822429906
$ws.Range("L17").Value = 22.33502538071
$ws.Range("M17").Value = 114.222222222222
$ws.Range("N17").Value = -16.608996539792
$ws.Range("F18").Value = 7
$ws.Range("G18").Value = 18
$ws.Range("H18").Value = -61.111111111111
$ws.Range("I18").Value = 157
$ws.Range("J18").Value = 142
$ws.Range("K18").Value = 10.56338028169
$ws.Range("L18").Value = -21.890547263681
$ws.Range("M18").Value = 65.263157894736
$ws.Range("N18").Value = -67.827868852459
$ws.Range("C19").Value = 4
$ws.Range("D19").Value = 4
$ws.Range("E19").Value = 0
$ws.Range("F19").Value = 21
$ws.Range("G19").Value = 28
$ws.Range("H19").Value = -25
$ws.Range("I19").Value = 373
$ws.Range("J19").Value = 435
$ws.Range("K19").Value = -14.252873563218
$ws.Range("L19").Value = -21.802935010482
$ws.Range("M19").Value = 45.136186770428
$ws.Range("N19").Value = -34.2151675485
$ws.Range("C20").Value = 2
$ws.Range("D20").Value = 3
$ws.Range("E20").Value = -33.333333333333
$ws.Range("G20").Value = 12
$ws.Range("H20").Value = -25
$ws.Range("I20").Value = 85
$ws.Range("J20").Value = 89
$ws.Range("K20").Value = -4.494382022471
$ws.Range("L20").Value = -4.494382022471
$ws.Range("M20").Value = 129.72972972973
$ws.Range("N20").Value = -79.21760391198
$ws.Range("D21").Value = 20
$ws.Range("E21").Value = 10
$ws.Range("F21").Value = 94
$ws.Range("G21").Value = 102
$ws.Range("H21").Value = -7.843137254901
$ws.Range("I21").Value = 1373
$ws.Range("J21").Value = 1319
$ws.Range("K21").Value = 4.094010614101
$ws.Range("L21").Value = 0.145878920495
$ws.Range("M21").Value = 52.895322939866
$ws.Range("N21").Value = -55.057283142389
$ws.Range("G22").Value = 1
$ws.Range("C23").Value = 10
$ws.Range("D23").Value = 6
$ws.Range("E23").Value = 66.666666666666
$ws.Range("F23").Value = 35
$ws.Range("G23").Value = 31
$ws.Range("H23").Value = 12.903225806451
$ws.Range("I23").Value = 410
$ws.Range("J23").Value = 425
$ws.Range("K23").Value = -3.529411764705
$ws.Range("L23").Value = 0.244498777506
$ws.Range("M23").Value = 38.983050847457
$ws.Range("C24").Value = 14
$ws.Range("E24").Value = -22.222222222222
$ws.Range("F24").Value = 63
$ws.Range("G24").Value = 71
$ws.Range("H24").Value = -11.267605633802
$ws.Range("I24").Value = 757
$ws.Range("J24").Value = 931
$ws.Range("K24").Value = -18.689581095596
$ws.Range("L24").Value = -17.806731813246
$ws.Range("M24").Value = 15.749235474006
$ws.Range("D25").Value = 1
$ws.Range("E25").Value = -100
$ws.Range("F25").Value = 13
$ws.Range("G25").Value = 12
$ws.Range("H25").Value = 8.333333333333
$ws.Range("J25").Value = 266
$ws.Range("K25").Value = -41.72932330827
$ws.Range("C26").Value = 1
$ws.Range("D26").Value = 7
$ws.Range("E26").Value = -85.714285714285
$ws.Range("F26").Value = 30
$ws.Range("G26").Value = 39
$ws.Range("H26").Value = -23.076923076923
$ws.Range("I26").Value = 732
$ws.Range("J26").Value = 659
$ws.Range("K26").Value = 11.077389984825
$ws.Range("L26").Value = 36.059479553903
$ws.Range("M26").Value = 9.417040358744
$ws.Range("F27").Value = 3
$ws.Range("I27").Value = 29
$ws.Range("K27").Value = 7.407407407407
$ws.Range("L27").Value = 7.407407407407
$ws.Range("F28").Value = 1
$ws.Range("H28").Value = -50
$ws.Range("C29").Value = 3
$ws.Range("F29").Value = 5
$ws.Range("H29").Value = 150
$ws.Range("I29").Value = 19
$ws.Range("K29").Value = 11.764705882352
$ws.Range("L29").Value = -26.923076923076
$ws.Range("M29").Value = -55.813953488372
$ws.Range("N29").Value = -80.41237113402
$ws.Range("C30").Value = 3
$ws.Range("F30").Value = 5
$ws.Range("H30").Value = 400
$ws.Range("I30").Value = 16
$ws.Range("K30").Value = 6.666666666666
$ws.Range("L30").Value = -15.78947368421
$ws.Range("M30").Value = -57.894736842105
$ws.Range("N30").Value = -82.417582417582
